$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Hours for the row 4 entry (was 4, now 2)
$ws.Range("C4").Value = 2

# Fill in the two new timetable rows (5 and 6), matching the author's
# original entry order so new shared-string entries line up
$ws.Range("A5").Value = "Wk [11] Sunday 27.5.18"
$ws.Range("D5").Value = "Designing Factory Class + Stats + Item"
$ws.Range("B5").Value = "1200 - 1530"
$ws.Range("C5").Value = 3.5

$ws.Range("A6").Value = "Wk [11] Sunday 27.5.18"
$ws.Range("B6").Value = "1700 - 1800"
$ws.Range("D6").Value = "Designing Stages"
$ws.Range("C6").Value = 1

# Move the active selection to C7, matching the author's next click
$ws.Range("C7").Select()
